# Fixed README.md stats and docx preparation for all Renaissance - JDK 17 -
# Shenandoah GC tests.
#
# The document is a single-column, single-row-per-value table. This edit
# updates a handful of summary-statistic cells, and collapses three rows
# that previously held a whole tab-separated stats line (count + 8 timing
# buckets + percentage) down to just their leading "count" value.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# --- Simple value replacements (rows 1-12, 1-indexed) ---
$t.Cell(1, 1).Range.Text  = "0M"
$t.Cell(2, 1).Range.Text  = "0M"
$t.Cell(3, 1).Range.Text  = "0M"
$t.Cell(4, 1).Range.Text  = "2390"
$t.Cell(5, 1).Range.Text  = "0.00002"
$t.Cell(6, 1).Range.Text  = "0.00967"
$t.Cell(7, 1).Range.Text  = "0.00028"
$t.Cell(8, 1).Range.Text  = "0.00042"
$t.Cell(9, 1).Range.Text  = "0.00024"
$t.Cell(10, 1).Range.Text = "0.00029"
$t.Cell(11, 1).Range.Text = "0.00040"
$t.Cell(12, 1).Range.Text = "0.81055"

# --- Collapse the three tab-separated rows down to a single value each ---
$t.Cell(44, 1).Range.Text = "99.93"
$t.Cell(45, 1).Range.Text = "0.81"
$t.Cell(46, 1).Range.Text = "1163"
